$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F2").Value = 1954
$ws1.Range("F7").Value = 1625
$ws1.Range("F8").Value = 22
$ws1.Range("F9").Value = 642
$ws1.Range("F14").Value = 224
$ws1.Range("F17").Value = 111
$ws1.Range("F19").Value = 3770
$ws1.Range("F21").Value = 17
$ws1.Range("F23").Value = 347
$ws1.Range("F24").Value = 702
$ws1.Range("F25").Value = 423
$ws1.Range("F28").Value = 1561
$ws1.Range("F30").Value = 150

$ws4.Range("F2").Value = 1954
$ws4.Range("F7").Value = 1625
$ws4.Range("F8").Value = 22
$ws4.Range("F9").Value = 642
$ws4.Range("F14").Value = 224
$ws4.Range("F17").Value = 111
$ws4.Range("F19").Value = 3770
$ws4.Range("F21").Value = 17
$ws4.Range("F23").Value = 347
$ws4.Range("F24").Value = 702
$ws4.Range("F25").Value = 423
$ws4.Range("F28").Value = 1561
$ws4.Range("F30").Value = 151
